# LVBR-85 kapitalisatie van notation en prefLabel
# For rows 7-23 (the "gebouw" concept rows), the notation (column J) is
# uppercased and the prefLabel (column L) is capitalized (first letter
# uppercase, remainder unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 7; $row -le 23; $row++) {
    $notationCell = $ws.Cells.Item($row, 10)   # column J = notation
    $prefLabelCell = $ws.Cells.Item($row, 12)  # column L = prefLabel

    $notation = $notationCell.Value2
    $prefLabel = $prefLabelCell.Value2

    if ($notation -ne $null) {
        $notationCell.Value = $notation.ToUpper().Replace(" ", "_")
    }

    if ($prefLabel -ne $null -and $prefLabel.Length -gt 0) {
        $capitalized = $prefLabel.Substring(0, 1).ToUpper() + $prefLabel.Substring(1)
        $prefLabelCell.Value = $capitalized
    }
}
